$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before current column D ("Tipo") to hold the new "MAE" column.
$ws.Range("D1").EntireColumn.Insert()

# Copy formatting from the neighboring header cell (C1) onto the new header cell (D1)
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header for the new column
$ws.Range("D1").Value = "MAE"

# New MAE values for rows 2-5
$ws.Range("D2").Value = 0.4321333824756292
$ws.Range("D3").Value = 0.2119198634755611
$ws.Range("D4").Value = 0.1361288253571671
$ws.Range("D5").Value = 0.1911874935925046

# Minor precision updates to existing MSE values (column B) for rows 4 and 5
$ws.Range("B4").Value = 0.04215534119371416
$ws.Range("B5").Value = 0.07796894984218661
